$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Harvard case classification: recompute the "*_old" app score columns
# (Ada_old, Avey_old, Buoy_old, K health_old, WebMD_old, doctor_MA_old,
# doctor_NJ_old, doctor_TH_old) and the average_doctor / average_doctor_old
# summary columns, and swap the average_doctor / average_doctor_old header
# labels (columns BP and BQ).

$ws.Range("BP1").Value = "average_doctor_old"
$ws.Range("BQ1").Value = "average_doctor"
$ws.Range("E4").Value = 0.613
$ws.Range("F4").Value = 0.049
$ws.Range("G4").Value = 0.221
$ws.Range("N4").Value = 0.633
$ws.Range("O4").Value = 0.093
$ws.Range("P4").Value = 0.306
$ws.Range("Q4").Value = 0.44
$ws.Range("R4").Value = 0.126
$ws.Range("S4").Value = 0.354
$ws.Range("W4").Value = 0.633
$ws.Range("X4").Value = 0.038
$ws.Range("Y4").Value = 0.194
$ws.Range("AI4").Value = 0.417
$ws.Range("AJ4").Value = 0.111
$ws.Range("AK4").Value = 0.333
$ws.Range("AU4").Value = 0.366
$ws.Range("AV4").Value = 0.007
$ws.Range("AW4").Value = 0.08400000000000001
$ws.Range("BA4").Value = 1.95
$ws.Range("BB4").Value = 0.06
$ws.Range("BC4").Value = 0.245
$ws.Range("BG4").Value = 0.733
$ws.Range("BH4").Value = 0.151
$ws.Range("BI4").Value = 0.389
$ws.Range("BM4").Value = 0.667
$ws.Range("BN4").Value = 0.078
$ws.Range("BO4").Value = 0.279
$ws.Range("BP4").Value = 0.65
$ws.Range("BQ4").Value = 0.798
$ws.Range("E5").Value = 0.713
$ws.Range("F5").Value = 0.062
$ws.Range("G5").Value = 0.249
$ws.Range("N5").Value = 0.6870000000000001
$ws.Range("O5").Value = 0.028
$ws.Range("P5").Value = 0.168
$ws.Range("Q5").Value = 0.29
$ws.Range("R5").Value = 0.048
$ws.Range("S5").Value = 0.218
$ws.Range("W5").Value = 0.423
$ws.Range("X5").Value = 0.08599999999999999
$ws.Range("Y5").Value = 0.293
$ws.Range("AI5").Value = 0.397
$ws.Range("AJ5").Value = 0.109
$ws.Range("AK5").Value = 0.331
$ws.Range("AU5").Value = 0.637
$ws.Range("AV5").Value = 0.057
$ws.Range("AW5").Value = 0.239
$ws.Range("BA5").Value = 1.106
$ws.Range("BB5").Value = 0.007
$ws.Range("BC5").Value = 0.083
$ws.Range("BG5").Value = 0.42
$ws.Range("BH5").Value = 0.057
$ws.Range("BI5").Value = 0.238
$ws.Range("BM5").Value = 0.323
$ws.Range("BN5").Value = 0.01
$ws.Range("BO5").Value = 0.102
$ws.Range("BP5").Value = 0.369
$ws.Range("BQ5").Value = 0.399
$ws.Range("E6").Value = 0.659
$ws.Range("N6").Value = 0.659
$ws.Range("Q6").Value = 0.35
$ws.Range("W6").Value = 0.507
$ws.Range("AI6").Value = 0.407
$ws.Range("AU6").Value = 0.465
$ws.Range("BA6").Value = 1.406
$ws.Range("BG6").Value = 0.534
$ws.Range("BM6").Value = 0.435
$ws.Range("BP6").Value = 0.469
$ws.Range("BQ6").Value = 0.53
$ws.Range("E7").Value = 0.6899999999999999
$ws.Range("N7").Value = 0.675
$ws.Range("Q7").Value = 0.311
$ws.Range("W7").Value = 0.453
$ws.Range("AI7").Value = 0.401
$ws.Range("AU7").Value = 0.555
$ws.Range("BA7").Value = 1.208
$ws.Range("BG7").Value = 0.459
$ws.Range("BM7").Value = 0.36
$ws.Range("BP7").Value = 0.403
$ws.Range("BQ7").Value = 0.443
$ws.Range("E8").Value = 0.8070000000000001
$ws.Range("F8").Value = 0.055
$ws.Range("G8").Value = 0.235
$ws.Range("N8").Value = 0.885
$ws.Range("O8").Value = 0.006
$ws.Range("P8").Value = 0.08
$ws.Range("Q8").Value = 0.319
$ws.Range("R8").Value = 0.118
$ws.Range("S8").Value = 0.343
$ws.Range("W8").Value = 0.734
$ws.Range("X8").Value = 0.002
$ws.Range("Y8").Value = 0.04
$ws.Range("AI8").Value = 0.511
$ws.Range("AJ8").Value = 0.157
$ws.Range("AK8").Value = 0.396
$ws.Range("AU8").Value = 0.654
$ws.Range("AV8").Value = 0.02
$ws.Range("AW8").Value = 0.141
$ws.Range("BA8").Value = 1.712
$ws.Range("BB8").Value = 0.048
$ws.Range("BC8").Value = 0.22
$ws.Range("BG8").Value = 0.611
$ws.Range("BH8").Value = 0.109
$ws.Range("BI8").Value = 0.33
$ws.Range("BM8").Value = 0.612
$ws.Range("BN8").Value = 0.066
$ws.Range("BO8").Value = 0.258
$ws.Range("BP8").Value = 0.571
$ws.Range("BQ8").Value = 0.616
$ws.Range("E9").Value = 0.8
$ws.Range("F9").Value = 0.16
$ws.Range("G9").Value = 0.4
$ws.Range("N9").Value = 1
$ws.Range("O9").Value = 0
$ws.Range("P9").Value = 0
$ws.Range("AI9").Value = 0.6
$ws.Range("BA9").Value = 1.8
$ws.Range("BM9").Value = 0.8
$ws.Range("BN9").Value = 0.16
$ws.Range("BO9").Value = 0.4
$ws.Range("BP9").Value = 0.6
$ws.Range("BQ9").Value = 0.667
$ws.Range("E10").Value = 0.8
$ws.Range("F10").Value = 0.16
$ws.Range("G10").Value = 0.4
$ws.Range("N10").Value = 1
$ws.Range("O10").Value = 0
$ws.Range("P10").Value = 0
$ws.Range("W10").Value = 1
$ws.Range("X10").Value = 0
$ws.Range("Y10").Value = 0
$ws.Range("AI10").Value = 0.6
$ws.Range("BA10").Value = 2.2
$ws.Range("BM10").Value = 0.8
$ws.Range("BN10").Value = 0.16
$ws.Range("BO10").Value = 0.4
$ws.Range("BP10").Value = 0.733
$ws.Range("BQ10").Value = 0.778
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 0
$ws.Range("N11").Value = 1
$ws.Range("O11").Value = 0
$ws.Range("P11").Value = 0
$ws.Range("W11").Value = 1
$ws.Range("X11").Value = 0
$ws.Range("Y11").Value = 0
$ws.Range("AI11").Value = 0.6
$ws.Range("AU11").Value = 0.8
$ws.Range("AV11").Value = 0.16
$ws.Range("AW11").Value = 0.4
$ws.Range("BA11").Value = 2.2
$ws.Range("BM11").Value = 0.8
$ws.Range("BN11").Value = 0.16
$ws.Range("BO11").Value = 0.4
$ws.Range("BP11").Value = 0.733
$ws.Range("BQ11").Value = 0.778
$ws.Range("E12").Value = 1.8
$ws.Range("F12").Value = 2.56
$ws.Range("G12").Value = 1.6
$ws.Range("W12").Value = 1.2
$ws.Range("X12").Value = 0.16
$ws.Range("Y12").Value = 0.4
$ws.Range("AU12").Value = 3.4
$ws.Range("AV12").Value = 5.04
$ws.Range("AW12").Value = 2.245
$ws.Range("BA12").Value = 3.833
$ws.Range("BB12").Value = 0.222
$ws.Range("BC12").Value = 0.471
$ws.Range("BP12").Value = 1.278
$ws.Range("BQ12").Value = 1.193
$ws.Range("BP13").Value = 0.65
$ws.Range("BQ13").Value = 0.538
